$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "Xbsh6UCb3l94ToInOCVi"
$ws.Range("E6").Value = "YpJsoRGg8G2DWU2PLZ78"

$ws.Range("D12").Value = "9JwUh0BdrG4KqCW7EIKQ"
$ws.Range("E12").Value = "lEhv6AtB5bMToBrwRe06"

$ws.Range("D18").Value = "CXMl6q4xEIJ2Lx51wGIB"
$ws.Range("E18").Value = "RaHCxD9RcyRjHk11IvTJ"

$ws.Range("E18").Select() | Out-Null
